$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'327.79"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'4.84%"
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value = "'4"
$ws.Range('G2').Style = 'Normal'
$ws.Range('D3').Value = "'39.74"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'5.85%"
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value = "'4"
$ws.Range('G3').Style = 'Normal'
$ws.Range('D4').Value = "'5.248"
$ws.Range('D4').Style = 'Normal'
$ws.Range('G4').Value = "'4"
$ws.Range('G4').Style = 'Normal'
$ws.Range('D5').Value = "'0.08099"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'2.35%"
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value = "'4"
$ws.Range('G5').Style = 'Normal'
$ws.Range('D6').Value = "'4.527"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'2.52%"
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value = "'4"
$ws.Range('G6').Style = 'Normal'
$ws.Range('D7').Value = "'8.634"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'4.37%"
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value = "'4"
$ws.Range('G7').Style = 'Normal'
$ws.Range('D8').Value = "'1.918"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'0.30%"
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value = "'4"
$ws.Range('G8').Style = 'Normal'
$ws.Range('E9').Value = "'-1.33%"
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value = "'4"
$ws.Range('G9').Style = 'Normal'
$ws.Range('D10').Value = "'0.9345"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'0.64%"
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value = "'4"
$ws.Range('G10').Style = 'Normal'
$ws.Range('D11').Value = "'0.1301"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'16.79%"
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value = "'4"
$ws.Range('G11').Style = 'Normal'
$ws.Range('D12').Value = "'0.1956"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'3.30%"
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value = "'4"
$ws.Range('G12').Style = 'Normal'
$ws.Range('D13').Value = "'0.09189"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'1.47%"
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value = "'4"
$ws.Range('G13').Style = 'Normal'
$ws.Range('D14').Value = "'0.03458"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'3.84%"
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value = "'4"
$ws.Range('G14').Style = 'Normal'
$ws.Range('D15').Value = "'0.09544"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'-0.57%"
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value = "'4"
$ws.Range('G15').Style = 'Normal'
$ws.Range('E16').Value = "'1.08%"
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value = "'4"
$ws.Range('G16').Style = 'Normal'
$ws.Range('D17').Value = "'0.04438"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'1.62%"
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value = "'4"
$ws.Range('G17').Style = 'Normal'
$ws.Range('D18').Value = "'0.005835"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'0.61%"
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value = "'4"
$ws.Range('G18').Style = 'Normal'
$ws.Range('G19').Value = "'4"
$ws.Range('G19').Style = 'Normal'
$ws.Range('E20').Value = "'3.71%"
$ws.Range('E20').Style = 'Normal'
$ws.Range('G20').Value = "'4"
$ws.Range('G20').Style = 'Normal'
$ws.Range('D21').Value = "'7.263"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'22.98%"
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value = "'4"
$ws.Range('G21').Style = 'Normal'
$ws.Range('E22').Value = "'1.31%"
$ws.Range('E22').Style = 'Normal'
$ws.Range('G22').Value = "'4"
$ws.Range('G22').Style = 'Normal'
$ws.Range('G23').Value = "'4"
$ws.Range('G23').Style = 'Normal'
$ws.Range('D24').Value = "'0.001223"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'-0.94%"
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value = "'4"
$ws.Range('G24').Style = 'Normal'
$ws.Range('D25').Value = "'0.004353"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'-6.19%"
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value = "'4"
$ws.Range('G25').Style = 'Normal'
$ws.Range('D26').Value = "'0.0001291"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'-5.12%"
$ws.Range('E26').Style = 'Normal'
$ws.Range('G26').Value = "'4"
$ws.Range('G26').Style = 'Normal'
$ws.Range('D27').Value = "'0.0003993"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'0.01%"
$ws.Range('E27').Style = 'Normal'
$ws.Range('G27').Value = "'4"
$ws.Range('G27').Style = 'Normal'
$ws.Range('G28').Value = "'4"
$ws.Range('G28').Style = 'Normal'
$ws.Range('G29').Value = "'4"
$ws.Range('G29').Style = 'Normal'
$ws.Range('G30').Value = "'4"
$ws.Range('G30').Style = 'Normal'
$ws.Range('G31').Value = "'4"
$ws.Range('G31').Style = 'Normal'
$ws.Range('G32').Value = "'4"
$ws.Range('G32').Style = 'Normal'
$ws.Range('G33').Value = "'4"
$ws.Range('G33').Style = 'Normal'
$ws.Range('G34').Value = "'4"
$ws.Range('G34').Style = 'Normal'
$ws.Range('G35').Value = "'4"
$ws.Range('G35').Style = 'Normal'
$ws.Range('G36').Value = "'4"
$ws.Range('G36').Style = 'Normal'
$ws.Range('G37').Value = "'4"
$ws.Range('G37').Style = 'Normal'
$ws.Range('G38').Value = "'4"
$ws.Range('G38').Style = 'Normal'
$ws.Range('D39').Value = "'0.02463"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'9.43%"
$ws.Range('E39').Style = 'Normal'
$ws.Range('G39').Value = "'4"
$ws.Range('G39').Style = 'Normal'
$ws.Range('D40').Value = "'0.05232"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'2.98%"
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value = "'4"
$ws.Range('G40').Style = 'Normal'
$ws.Range('D41').Value = "'0.007714"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'3.50%"
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value = "'4"
$ws.Range('G41').Style = 'Normal'
$ws.Range('D42').Value = "'0.1433"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'5.76%"
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value = "'4"
$ws.Range('G42').Style = 'Normal'
$ws.Range('D43').Value = "'0.008677"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'-3.74%"
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value = "'4"
$ws.Range('G43').Style = 'Normal'
$ws.Range('D44').Value = "'0.002111"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'-0.91%"
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value = "'4"
$ws.Range('G44').Style = 'Normal'
$ws.Range('D45').Value = "'0.008175"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'-5.32%"
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value = "'4"
$ws.Range('G45').Style = 'Normal'
$ws.Range('D46').Value = "'0.00006666"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'-0.37%"
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value = "'4"
$ws.Range('G46').Style = 'Normal'
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'0.02%"
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value = "'4"
$ws.Range('G47').Style = 'Normal'
$ws.Range('D48').Value = "'0.002853"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'-13.10%"
$ws.Range('E48').Style = 'Normal'
$ws.Range('G48').Value = "'4"
$ws.Range('G48').Style = 'Normal'
$ws.Range('E49').Value = "'148.14%"
$ws.Range('E49').Style = 'Normal'
$ws.Range('G49').Value = "'4"
$ws.Range('G49').Style = 'Normal'
$ws.Range('D50').Value = "'0.00002102"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'0.02%"
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value = "'4"
$ws.Range('G50').Style = 'Normal'
$ws.Range('D51').Value = "'0.0002001"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'0.02%"
$ws.Range('E51').Style = 'Normal'
$ws.Range('G51').Value = "'4"
$ws.Range('G51').Style = 'Normal'
